$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 55, "Tax/BaseScale")
# onto the two new rows so they inherit the same font/fill treatment.
$ws.Range("B55:D55").Copy()
$ws.Range("B56:D57").PasteSpecial(-4122)
$ws.Range("B56:D57").Borders.LineStyle = -4142

# New parameter: MFTC_WEP_scaling
$ws.Range("B56").Value = "MFTC_WEP_scaling"
$ws.Range("C56").Value = "1"
$ws.Range("D56").Value = "How should the Winter Energy Payment be scaled? Average week = 1, Winter week = 12/5, Summer week = 0"

# New parameter: WFF_or_Benefit
$ws.Range("B57").Value = "WFF_or_Benefit"
$ws.Range("C57").Value = "Max"
$ws.Range("D57").Value = 'What work decision should we assume? Go off-benefit and receive IWTC = "WFF", stay on-benefit = "Benefit", or whichever gives a higher net income = "Max"'
